# CV_Frontend_Yakima_EN.docx edit script
# Implements the changes described by the commit:
#  1. Remove the stray "_GoBack" bookmark that used to sit between the
#     "Home" and "Page" runs.
#  2. Merge the 3-run hyperlink text "http://www.yxey" + "e" + ".com/"
#     into a single run "http://www.yxeye.com/" (keeping the Hyperlink
#     character style).
#  3. Insert "北京" before the existing "无线" run (so the company name
#     reads "...@ 北京无线天利有限公司上海分公司"), producing a brand new
#     run for "北京" while keeping "无线" and "天利有限公司上海分公司" as
#     their own separate runs (matching the source formatting).
#  4. Split " project outsourced by China Pacific Insurance using
#     ReactNative+Redux;" into " project outso" / "urced by China
#     Pacific Insurance using ReactNative+Redux;" with a new "_GoBack"
#     bookmark sitting at the split point (this is where Word's cursor
#     was left after editing, per the diff).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Delete the old "_GoBack" bookmark (originally right after "Home").
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Collapse the 3 hyperlink runs into one clean run, preserving the
#    "Hyperlink" character style (aff5).
# ---------------------------------------------------------------------
$find = $d.Content
$null = $find.Find.Execute("http://www.yxey" + "e" + ".com/", $false, $false, $false, $false, $false, $true, 1, $false, "http://www.yxeye.com/", 2)

$styleRng = $d.Content
$null = $styleRng.Find.Execute("http://www.yxeye.com/", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$styleRng.Style = "aff5"

# ---------------------------------------------------------------------
# 3) "北京" + "无线" + "天利有限公司上海分公司" as three separate runs.
# ---------------------------------------------------------------------
$wx = $d.Content
$null = $wx.Find.Execute("无线", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$wxStart = $wx.Start
$wxEnd = $wx.End

# Rename the existing "无线" run's text to "北京" (keeps that run, and
# its formatting, intact).
$renameRng = $d.Range($wxStart, $wxEnd)
$renameRng.Text = "北京"

# Insert a fresh "无线" right after it (inherits the same run's
# formatting because it is typed inside that run).
$insertPoint = $d.Range($wxStart + 2, $wxStart + 2)
$insertPoint.InsertBefore("无线")

# The single text edit above merges "北京无线" and the following
# "天利有限公司上海分公司" run into one run (identical formatting), so
# force the two splits back apart with a formatting no-op (set+unset a
# property) which re-splits runs without re-merging the whole
# paragraph.
$beijingRng = $d.Range($wxStart, $wxStart + 2)
$beijingRng.Bold = 1
$beijingRng.Bold = 0

$wuxianRng = $d.Range($wxStart + 2, $wxStart + 4)
$wuxianRng.Bold = 1
$wuxianRng.Bold = 0

# ---------------------------------------------------------------------
# 4) Split "... project outsourced by China Pacific Insurance ..." at
#    "outso|urced" and drop a new "_GoBack" bookmark at the split.
# ---------------------------------------------------------------------
$outso = $d.Content
$null = $outso.Find.Execute("outso", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $outso.End

$bmRng = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRng)

# Re-assert the trailing run's text (self-replace) so it normalizes to
# a plain <w:t> (no stray xml:space="preserve") like the reference.
$tail = $d.Content
$null = $tail.Find.Execute("urced by China Pacific Insurance using ReactNative+Redux;", $false, $false, $false, $false, $false, $true, 1, $false, "urced by China Pacific Insurance using ReactNative+Redux;", 2)

Write-Output "done"
